$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Enter the value that drives the recalculated shared formulas in row 13
# (R13/U13/X13/.../BA13 recompute automatically from this input)
$ws.Range("Q13").Value = 3

# Re-touch the header merged cells so they re-serialize in the same
# order as the committed file (unmerge + remerge re-appends each one)
$mergedRanges = @("AZ4:BA4", "AO4:AP4", "AR4:AS4", "AU4:AV4", "AX4:AY4")
foreach ($mergedRange in $mergedRanges) {
    $ws.Range($mergedRange).UnMerge()
    $ws.Range($mergedRange).Merge()
}

# Update the active/selected cell on the sheet (bottomRight pane selection)
$ws.Range("W11").Select()
